$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the single "Interlocuteur" header with the new set of column headers
$ws.Range("A1").Value = "email"
$ws.Range("B1").Value = "nom"
$ws.Range("C1").Value = "prénom"
$ws.Range("D1").Value = "civilite"
$ws.Range("E1").Value = "tel"
$ws.Range("F1").Value = "adresse"

# Update selection to match the new active cell (F1)
$null = $ws.Range("F1").Select()
